# fb-surveyState.xlsx update: "data up to 5th"
# Extends the per-state daily survey-response table (Sheet1) through 05 08 2020,
# adds the new "Puerto Rico" (column AR) figures that start appearing with the
# 21 07 2020 - 28 07 2020 week, and re-syncs seven previously-estimated values in
# the 27 07 2020 row (178) now that the underlying survey wave finished processing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Puerto Rico" column (AR) gains its first values for rows 153-160
#    (the 13 07 2020 - 20 07 2020 week); AR stayed blank for every earlier row.
# ---------------------------------------------------------------------------
$puertoRicoByRow = [ordered]@{
    153 = 0.2155172
    154 = 0.2403846
    155 = 0.245098
    156 = 0.2252252
    157 = 0.2173913
    158 = 0.2272727
    159 = 0
    160 = 0
}
foreach ($row in $puertoRicoByRow.Keys) {
    $ws.Range("AR" + $row).Value = $puertoRicoByRow[$row]
}

# ---------------------------------------------------------------------------
# 2) Row 178 (27 07 2020) - seven columns get corrected/refined values.
# ---------------------------------------------------------------------------
$row178Updates = [ordered]@{
    "M" = 0.9341245
    "U" = 0.4794321
    "W" = 0.2445872
    "AO" = 0.8842949
    "AV" = 0.7224034
    "AW" = 0.9666351
    "BB" = 0.3523222
}
foreach ($col in $row178Updates.Keys) {
    $ws.Range($col + "178").Value = $row178Updates[$col]
}

# ---------------------------------------------------------------------------
# 3) New daily rows 179-188, one per date from 27 07 2020 to 05 08 2020.
#    Row 179 (27 07 2020) already had its date label; it only needed the rest
#    of its survey figures. Rows 180-186 are brand-new data rows. Rows 187-188
#    (04 08 2020 and 05 08 2020) only carry the date label so far.
# ---------------------------------------------------------------------------
$dateLabels = @(
    "28 07 2020", "29 07 2020", "30 07 2020", "31 07 2020", "01 08 2020",
    "02 08 2020", "03 08 2020", "04 08 2020", "05 08 2020"
)
for ($i = 0; $i -lt $dateLabels.Count; $i++) {
    $ws.Range("A" + (180 + $i)).Value = $dateLabels[$i]
}

# Columns that carry data in every survey row (E, N, AC, AR and AZ are always
# blank in this sheet, e.g. unused/retired state codes).
$dataColumns = @(
    "B", "C", "D", "F", "G", "H", "I", "J", "K", "L",
    "M", "O", "P", "Q", "R", "S", "T", "U", "V", "W",
    "X", "Y", "Z", "AA", "AB", "AD", "AE", "AF", "AG", "AH",
    "AI", "AJ", "AK", "AL", "AM", "AN", "AO", "AP", "AQ", "AS",
    "AT", "AU", "AV", "AW", "AX", "AY", "BA", "BB", "BC", "BD",
    "BE"
)

$rowValues = [ordered]@{
    179 = @(0.5248316, 1.212313, 0.8538848999999999, 0.7013397, 0.6103184, 0.4242801, 0.2509629, 0.286212, 0.5529746, 0.8730049, 0.9511316, 0.3784477, 0.8039177, 0.8396976, 0.4702863, 0.4288657, 0.5213249, 0.4121954, 1.1911954, 0.2371935, 0.3572759, 0.2259051, 0.4646143, 0.3744387, 0.6519963, 1.1971962, 0.520822, 0.4952999, 0.8563452, 0.7149718, 0.2458613, 0.2944996, 0.5572198, 0.5250286, 0.3309999, 0.5566467, 0.8032484, 0.3961243, 0.3865759, 0.3897981, 0.9437637, 0.431679, 0.7412112, 0.9249517, 0.6245039999999999, 0.4340576, 0.3722836, 0.3380064, 0.4527238, 0.5393454, 0.7196131)
    180 = @(0.3852219, 1.1438661, 0.8384449, 0.6444778, 0.6354145, 0.3995012, 0.2818689, 0.288558, 0.5532305, 0.8587176, 0.9803891, 0.2930057, 0.7782866000000001, 0.8125331, 0.4315586, 0.4717795, 0.5406936, 0.4612975, 1.1717163, 0.1985704, 0.3636416, 0.2052703, 0.4946733, 0.4402342, 0.6325896, 1.1754554, 0.4812375, 0.5120221, 1.0294477, 0.6722232, 0.275387, 0.3036089, 0.563046, 0.4900312, 0.3119919, 0.5461734, 0.9184307, 0.3698854, 0.3291252, 0.3879592, 0.8829156, 0.4347915, 0.7301528, 0.8864405, 0.5597618, 0.4106016, 0.3801402, 0.3602823, 0.4490463, 0.4775579, 0.8478426)
    181 = @(0.5402273, 1.0903966, 0.7844236999999999, 0.6648596, 0.6432824, 0.4009179, 0.2197542, 0.1139818, 0.4117909, 0.8617093, 1.0331518, 0.2680251, 0.7276541, 0.8919244, 0.4616245, 0.4825277, 0.45925, 0.4499493, 1.1165638, 0.2054785, 0.3437636, 0.1338878, 0.5191832, 0.4352785, 0.6735883, 1.189163, 0.5393945999999999, 0.5343337, 1.0360806, 0.6113571, 0.282106, 0.2888616, 0.4959295, 0.618947, 0.3218038, 0.5130676, 0.8858655, 0.3767333, 0.3095264, 0.3898234, 0.8669355, 0.5540519, 0.7498977999999999, 0.8902594, 0.5577701, 0.3811209, 0.3213503, 0.3479924, 0.5186669, 0.5486291, 0.8803368)
    182 = @(0.647823, 1.0539639, 0.7204587, 0.6121416, 0.6344083, 0.4023165, 0.188324, 0.1147959, 0.3444314, 0.8676807, 1.07292, 0.2414923, 0.7027792, 0.7457607000000001, 0.4579933, 0.5132407, 0.4218282, 0.54832, 1.047941, 0.2044087, 0.3394891, 0.1361392, 0.5166516, 0.5155669000000001, 0.5835835, 1.1205245, 0.6256931, 0.5543896, 1.0013839, 0.5865977, 0.2685616, 0.3089803, 0.483274, 0.6831631, 0.2981458, 0.5337483, 0.8798565, 0.4350771, 0.3287845, 0.383964, 0.773935, 0.49758, 0.7289563999999999, 0.8816857, 0.5726809, 0.3740717, 0.4152631, 0.3471397, 0.5014328, 0.550752, 0.7411326)
    183 = @(0.6653439, 1.0455662, 0.7400926, 0.5941626, 0.6399652, 0.4321611, 0.1862579, 0.1133121, 0.3673204, 0.8392395, 1.0903209, 0.2133211, 0.6377775, 0.7782686, 0.4636993, 0.5391241, 0.4247705, 0.6401936, 1.0615493, 0.2159326, 0.3196161, 0.1271579, 0.4962813, 0.5634726, 0.5157739, 1.2068702, 0.585122, 0.5622837000000001, 1.0050089, 0.5968405, 0.3492996, 0.3434478, 0.4611194, 0.6423428, 0.2757025, 0.5629174, 0.7903992, 0.5092869, 0.3670359, 0.3735477, 0.8289143, 0.5563715, 0.8296396, 0.8607769, 0.5205376, 0.3835031, 0.5114933, 0.3432314, 0.5063105, 0.4987776, 0.5429049)
    184 = @(0.7337373, 1.0691739, 0.8235914, 0.5135164, 0.6509644999999999, 0.4041627, 0.1993346, 0.1315061, 0.2601787, 0.8419316999999999, 1.0990987, 0.2692539, 0.601369, 0.8525862, 0.4522455, 0.6161093, 0.3465708, 0.7039264, 0.9794647, 0.2021601, 0.3082337, 0.0564831, 0.4813592, 0.6096499, 0.4883653, 1.1372292, 0.6592988, 0.5565145, 0.9566281, 0.6571862000000001, 0.3363014, 0.3415756, 0.373749, 0.7218037, 0.2725357, 0.5641993, 0.8283014, 0.4837814, 0.3881674, 0.3849761, 0.8396607, 0.5117436, 0.877451, 0.8609422, 0.5434701, 0.3692099, 0.4366743, 0.301121, 0.4445167, 0.5615816, 0.4663869)
    185 = @(0.7071153, 1.0737752, 0.7833976, 0.5236605, 0.6329932, 0.4169433, 0.2208385, 0.2472103, 0.4275793, 0.8311790999999999, 1.1077173, 0.2451865, 0.6152466, 0.7756389, 0.4527189, 0.6685207, 0.3746816, 0.7220658, 1.0131189, 0.232009, 0.319815, 0.081543, 0.5147445, 0.5289835000000001, 0.5463573, 1.1786656, 0.6197173, 0.5613074, 1.1299567, 0.6818988, 0.3320492, 0.2788197, 0.4277172, 0.7542908, 0.3140345, 0.5514277, 0.8353341, 0.4504935, 0.3695224, 0.3680153, 0.85522, 0.6758507, 0.8749504, 0.8850225, 0.5624533, 0.3688266, 0.3916166, 0.2835486, 0.4767237, 0.5463871, 0.4379936)
    186 = @(0.7807181, 1.0462526, 0.8035135, 0.4999616, 0.6592538, 0.4104711, 0.2133732, 0.2573238, 0.3937372, 0.7886992, 1.0524196, 0.2070064, 0.6401834, 0.9024922, 0.4252655, 0.7135096, 0.4358924, 0.766464, 0.9307893, 0.2397688, 0.3155382, 0.1128086, 0.4349477, 0.5576414, 0.5216417, 1.2224889, 0.6884428, 0.5950555, 1.0107168, 0.6040753, 0.2592929, 0.2625301, 0.3788708, 0.8592483, 0.2985399, 0.517974, 0.8302994, 0.3725367, 0.3837196, 0.3520063, 0.8879481, 0.7488973, 0.8239841, 0.895998, 0.5633666000000001, 0.3776351, 0.3877404, 0.3001878, 0.4494349, 0.4668688, 0.4889031)
}

foreach ($row in $rowValues.Keys) {
    $values = $rowValues[$row]
    for ($i = 0; $i -lt $dataColumns.Count; $i++) {
        $ws.Range($dataColumns[$i] + $row).Value = $values[$i]
    }
}

Write-Host "fb-surveyState: extended Sheet1 through 05 08 2020 (rows 153-188)."